$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("June")

$ws.Range("B2").Value = 1367
$ws.Range("C2").Value = 1151
$ws.Range("D2").Value = 216
$ws.Range("E2").Value = "We borrowerd more than we lent"
$ws.Range("G2").Value = "1.19 : 1"

$ws.Range("B3").Value = 648
$ws.Range("C3").Value = 370
$ws.Range("D3").Value = 278
$ws.Range("E3").Value = "We borrowerd more than we lent"
$ws.Range("G3").Value = "1.75 : 1"

$ws.Range("B4").Value = 1120
$ws.Range("C4").Value = 1014
$ws.Range("D4").Value = 106
$ws.Range("E4").Value = "We borrowerd more than we lent"
$ws.Range("G4").Value = "1.10 : 1"

$ws.Range("B5").Value = 44
$ws.Range("C5").Value = 171
$ws.Range("D5").Value = -127
$ws.Range("F5").Value = "We lent more than we borrowed"
$ws.Range("G5").Value = "0.26 : 1"

$ws.Range("B6").Value = 1008
$ws.Range("C6").Value = 1389
$ws.Range("D6").Value = -381
$ws.Range("F6").Value = "We lent more than we borrowed"
$ws.Range("G6").Value = "0.73 : 1"

$ws.Range("B7").Value = 259
$ws.Range("C7").Value = 153
$ws.Range("D7").Value = 106
$ws.Range("E7").Value = "We borrowerd more than we lent"
$ws.Range("G7").Value = "1.69 : 1"

$ws.Range("B8").Value = 177
$ws.Range("C8").Value = 205
$ws.Range("D8").Value = -28
$ws.Range("F8").Value = "We lent more than we borrowed"
$ws.Range("G8").Value = "0.86 : 1"

$ws.Range("B9").Value = 45
$ws.Range("C9").Value = 69
$ws.Range("D9").Value = -24
$ws.Range("F9").Value = "We lent more than we borrowed"
$ws.Range("G9").Value = "0.65 : 1"

$ws.Range("B10").Value = 7
$ws.Range("C10").Value = 38
$ws.Range("D10").Value = -31
$ws.Range("F10").Value = "We lent more than we borrowed"
$ws.Range("G10").Value = "0.18 : 1"

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 36
$ws.Range("C12").Value = 16
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = "We borrowerd more than we lent"
$ws.Range("G12").Value = "2.25 : 1"

$ws.Range("B13").Value = 141
$ws.Range("C13").Value = 86
$ws.Range("D13").Value = 55
$ws.Range("E13").Value = "We borrowerd more than we lent"
$ws.Range("G13").Value = "1.64 : 1"

$ws.Range("B14").Value = 97
$ws.Range("C14").Value = 242
$ws.Range("D14").Value = -145
$ws.Range("F14").Value = "We lent more than we borrowed"
$ws.Range("G14").Value = "0.40 : 1"

$ws.Range("B15").Value = 82
$ws.Range("C15").Value = 134
$ws.Range("D15").Value = -52
$ws.Range("F15").Value = "We lent more than we borrowed"
$ws.Range("G15").Value = "0.61 : 1"

$ws.Range("B16").Value = 65
$ws.Range("C16").Value = 147
$ws.Range("D16").Value = -82
$ws.Range("F16").Value = "We lent more than we borrowed"
$ws.Range("G16").Value = "0.44 : 1"

$ws.Range("B17").Value = 573
$ws.Range("C17").Value = 466
$ws.Range("D17").Value = 107
$ws.Range("E17").Value = "We borrowerd more than we lent"
$ws.Range("G17").Value = "1.23 : 1"

$ws.Range("B18").Value = 71
$ws.Range("C18").Value = 126
$ws.Range("D18").Value = -55
$ws.Range("F18").Value = "We lent more than we borrowed"
$ws.Range("G18").Value = "0.56 : 1"

$ws.Range("B19").Value = 482
$ws.Range("C19").Value = 359
$ws.Range("D19").Value = 123
$ws.Range("E19").Value = "We borrowerd more than we lent"
$ws.Range("G19").Value = "1.34 : 1"

$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 59
$ws.Range("D20").Value = -58
$ws.Range("F20").Value = "We lent more than we borrowed"
$ws.Range("G20").Value = "0.02 : 1"

$ws.Range("B21").Value = 640
$ws.Range("C21").Value = 357
$ws.Range("D21").Value = 283
$ws.Range("E21").Value = "We borrowerd more than we lent"
$ws.Range("G21").Value = "1.79 : 1"

$ws.Range("B22").Value = 23
$ws.Range("C22").Value = 82
$ws.Range("D22").Value = -59
$ws.Range("F22").Value = "We lent more than we borrowed"
$ws.Range("G22").Value = "0.28 : 1"

$ws.Range("B23").Value = 668
$ws.Range("C23").Value = 373
$ws.Range("D23").Value = 295
$ws.Range("E23").Value = "We borrowerd more than we lent"
$ws.Range("G23").Value = "1.79 : 1"

$ws.Range("B24").Value = 1338
$ws.Range("C24").Value = 1274
$ws.Range("D24").Value = 64
$ws.Range("E24").Value = "We borrowerd more than we lent"
$ws.Range("G24").Value = "1.05 : 1"

$ws.Range("B25").Value = 182
$ws.Range("C25").Value = 370
$ws.Range("D25").Value = -188
$ws.Range("F25").Value = "We lent more than we borrowed"
$ws.Range("G25").Value = "0.49 : 1"

$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = 277
$ws.Range("C27").Value = 189
$ws.Range("D27").Value = 88
$ws.Range("E27").Value = "We borrowerd more than we lent"
$ws.Range("G27").Value = "1.47 : 1"

$ws.Range("B28").Value = 61
$ws.Range("C28").Value = 96
$ws.Range("D28").Value = -35
$ws.Range("F28").Value = "We lent more than we borrowed"
$ws.Range("G28").Value = "0.64 : 1"

$ws.Range("B29").Value = 394
$ws.Range("C29").Value = 425
$ws.Range("D29").Value = -31
$ws.Range("F29").Value = "We lent more than we borrowed"
$ws.Range("G29").Value = "0.93 : 1"

$ws.Range("B30").Value = 44
$ws.Range("C30").Value = 26
$ws.Range("D30").Value = 18
$ws.Range("E30").Value = "We borrowerd more than we lent"
$ws.Range("G30").Value = "1.69 : 1"

$ws.Range("B31").Value = 32
$ws.Range("C31").Value = 258
$ws.Range("D31").Value = -226
$ws.Range("F31").Value = "We lent more than we borrowed"
$ws.Range("G31").Value = "0.12 : 1"

$ws.Range("B32").Value = 408
$ws.Range("C32").Value = 611
$ws.Range("D32").Value = -203
$ws.Range("F32").Value = "We lent more than we borrowed"
$ws.Range("G32").Value = "0.67 : 1"

$ws.Range("B33").Value = 385
$ws.Range("C33").Value = 565
$ws.Range("D33").Value = -180
$ws.Range("F33").Value = "We lent more than we borrowed"
$ws.Range("G33").Value = "0.68 : 1"

$ws.Range("B34").Value = 186
$ws.Range("C34").Value = 74
$ws.Range("D34").Value = 112
$ws.Range("E34").Value = "We borrowerd more than we lent"
$ws.Range("G34").Value = "2.51 : 1"

$ws.Range("B35").Value = 762
$ws.Range("C35").Value = 1002
$ws.Range("D35").Value = -240
$ws.Range("F35").Value = "We lent more than we borrowed"
$ws.Range("G35").Value = "0.76 : 1"

$ws.Range("B36").Value = 223
$ws.Range("C36").Value = 453
$ws.Range("D36").Value = -230
$ws.Range("F36").Value = "We lent more than we borrowed"
$ws.Range("G36").Value = "0.49 : 1"

$ws.Range("B37").Value = 494
$ws.Range("C37").Value = 265
$ws.Range("D37").Value = 229
$ws.Range("E37").Value = "We borrowerd more than we lent"
$ws.Range("G37").Value = "1.86 : 1"

$ws.Range("B38").Value = 35
$ws.Range("C38").Value = 184
$ws.Range("D38").Value = -149
$ws.Range("F38").Value = "We lent more than we borrowed"
$ws.Range("G38").Value = "0.19 : 1"

$ws.Range("B39").Value = 0
$ws.Range("C39").Value = 10
$ws.Range("D39").Value = -10
$ws.Range("F39").Value = "We lent more than we borrowed"
$ws.Range("G39").Value = "0.00 : 1"

$ws.Range("B40").Value = 0
$ws.Range("C40").Value = 9
$ws.Range("D40").Value = -9
$ws.Range("F40").Value = "We lent more than we borrowed"
$ws.Range("G40").Value = "0.00 : 1"

$ws.Range("B41").Value = 0
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = -2
$ws.Range("F41").Value = "We lent more than we borrowed"
$ws.Range("G41").Value = "0.00 : 1"

$ws.Range("B42").Value = 0
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = -2
$ws.Range("F42").Value = "We lent more than we borrowed"
$ws.Range("G42").Value = "0.00 : 1"

$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0

$ws.Range("B44").Value = 54
$ws.Range("C44").Value = 73
$ws.Range("D44").Value = -19
$ws.Range("F44").Value = "We lent more than we borrowed"
$ws.Range("G44").Value = "0.74 : 1"

$ws.Range("B45").Value = 91
$ws.Range("C45").Value = 189
$ws.Range("D45").Value = -98
$ws.Range("F45").Value = "We lent more than we borrowed"
$ws.Range("G45").Value = "0.48 : 1"

$ws.Range("B46").Value = 445
$ws.Range("C46").Value = 524
$ws.Range("D46").Value = -79
$ws.Range("F46").Value = "We lent more than we borrowed"
$ws.Range("G46").Value = "0.85 : 1"

$ws.Range("B47").Value = 1034
$ws.Range("C47").Value = 525
$ws.Range("D47").Value = 509
$ws.Range("E47").Value = "We borrowerd more than we lent"
$ws.Range("G47").Value = "1.97 : 1"

$ws.Range("B48").Value = 215
$ws.Range("C48").Value = 497
$ws.Range("D48").Value = -282
$ws.Range("F48").Value = "We lent more than we borrowed"
$ws.Range("G48").Value = "0.43 : 1"

$ws.Range("B49").Value = 504
$ws.Range("C49").Value = 276
$ws.Range("D49").Value = 228
$ws.Range("E49").Value = "We borrowerd more than we lent"
$ws.Range("G49").Value = "1.83 : 1"

$ws.Range("B50").Value = 969
$ws.Range("C50").Value = 608
$ws.Range("D50").Value = 361
$ws.Range("E50").Value = "We borrowerd more than we lent"
$ws.Range("G50").Value = "1.59 : 1"

$ws.Range("B51").Value = 203
$ws.Range("C51").Value = 161
$ws.Range("D51").Value = 42
$ws.Range("E51").Value = "We borrowerd more than we lent"
$ws.Range("G51").Value = "1.26 : 1"

$ws.Range("B52").Value = 320
$ws.Range("C52").Value = 407
$ws.Range("D52").Value = -87
$ws.Range("F52").Value = "We lent more than we borrowed"
$ws.Range("G52").Value = "0.79 : 1"

$ws.Range("B53").Value = 103
$ws.Range("C53").Value = 237
$ws.Range("D53").Value = -134
$ws.Range("F53").Value = "We lent more than we borrowed"
$ws.Range("G53").Value = "0.43 : 1"

$ws.Range("B54").Value = 28
$ws.Range("C54").Value = 213
$ws.Range("D54").Value = -185
$ws.Range("F54").Value = "We lent more than we borrowed"
$ws.Range("G54").Value = "0.13 : 1"

$ws.Range("B55").Value = 363
$ws.Range("C55").Value = 172
$ws.Range("D55").Value = 191
$ws.Range("E55").Value = "We borrowerd more than we lent"
$ws.Range("G55").Value = "2.11 : 1"

$ws.Activate()
